$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 7138
$ws1.Range("F4").Value = 3483
$ws1.Range("F14").Value = 114
$ws1.Range("F21").Value = 402
$ws1.Range("F24").Value = 1615
$ws1.Range("F27").Value = 2956
$ws1.Range("F28").Value = 2134
$ws1.Range("F35").Value = 4162
$ws1.Range("F40").Value = 744
$ws1.Range("F41").Value = 150

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 29
$ws2.Range("F15").Value = 551

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 7138
$ws4.Range("F6").Value = 3483
$ws4.Range("F15").Value = 114
$ws4.Range("F24").Value = 402
$ws4.Range("F26").Value = 1615
$ws4.Range("F29").Value = 2956
$ws4.Range("F30").Value = 2134
$ws4.Range("F36").Value = 4162
$ws4.Range("F42").Value = 744
